$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Give the merged-range "filler" cells (C1, D1) a top+bottom border, and make the
# trailing cell of the merge (D1) additionally carry a right border, mirroring the
# border treatment already used elsewhere in the sheet.
$ws1.Range("C1").Borders.Item(8).LineStyle = 1
$ws1.Range("C1").Borders.Item(9).LineStyle = 1

$ws1.Range("D1").Borders.Item(8).LineStyle = 1
$ws1.Range("D1").Borders.Item(9).LineStyle = 1
$ws1.Range("D1").Borders.Item(10).LineStyle = 1

# Anonymize the "fedcore" column header.
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison -------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Borders.Item(8).LineStyle = 1
$ws2.Range("C1").Borders.Item(9).LineStyle = 1

$ws2.Range("D1").Borders.Item(8).LineStyle = 1
$ws2.Range("D1").Borders.Item(9).LineStyle = 1
$ws2.Range("D1").Borders.Item(10).LineStyle = 1

$ws2.Range("F1").Borders.Item(8).LineStyle = 1
$ws2.Range("F1").Borders.Item(9).LineStyle = 1

$ws2.Range("G1").Borders.Item(8).LineStyle = 1
$ws2.Range("G1").Borders.Item(9).LineStyle = 1
$ws2.Range("G1").Borders.Item(10).LineStyle = 1

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell left behind in the model_size row.
$ws2.Range("G5").ClearContents()
